$p = $ppt.ActivePresentation

# Slide 1: Title "Header" + " " + "with" + " " + [Courier]"inline code"
#   -> merge the first four runs into a single "Header with " run,
#      keep the Courier-formatted "inline code" run untouched.
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Characters(1, 12).Text = "Header with "

# Slide 2: Title "Syntax" + " " + "highlighting" -> single run "Syntax highlighting"
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Characters(1, 19).Text = "Syntax highlighting"

# Slide 3: Title "Two" + " " + "column" + " " + "slide" -> single run "Two column slide"
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Characters(1, 17).Text = "Two column slide"
